# berlin52 TSP "long" results sheet — refresh with a new GA run.
# The genetic algorithm ran again ("GA FOREVAAA AND EVAAA"); this produced
# a new best tour over the same 52 cities with a (worse, in this case)
# total tour length of 9663.8544921875, replacing the previous
# 7469.65771484375 result. Row 1 holds the city count, row 2 the summary
# (best length in E2), row 3 mirrors the best length into D3, and rows
# 4-55 hold the new tour order (col A = city id) together with each
# city's coordinates (cols C/D) and a constant flag (col E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Summary / header rows: new best tour length.
$ws.Range("E2").Value = 9663.8544921875
$ws.Range("D3").Value = 9663.8544921875

# New tour: city id (A), unused flag (B), x (C), y (D), constant (E).
$tour = @(
    @(30.0, 0.0, 410.0, 250.0, 1.0),
    @(29.0, 0.0, 660.0, 180.0, 1.0),
    @(47.0, 0.0, 1170.0, 65.0, 1.0),
    @(26.0, 0.0, 1215.0, 245.0, 1.0),
    @(14.0, 0.0, 1530.0, 5.0, 1.0),
    @(52.0, 0.0, 1740.0, 245.0, 1.0),
    @(13.0, 0.0, 1465.0, 200.0, 1.0),
    @(27.0, 0.0, 1320.0, 315.0, 1.0),
    @(28.0, 0.0, 1250.0, 400.0, 1.0),
    @(11.0, 0.0, 1605.0, 620.0, 1.0),
    @(51.0, 0.0, 1340.0, 725.0, 1.0),
    @(12.0, 0.0, 1220.0, 580.0, 1.0),
    @(16.0, 0.0, 725.0, 370.0, 1.0),
    @(46.0, 0.0, 830.0, 485.0, 1.0),
    @(44.0, 0.0, 700.0, 500.0, 1.0),
    @(50.0, 0.0, 595.0, 360.0, 1.0),
    @(20.0, 0.0, 560.0, 365.0, 1.0),
    @(23.0, 0.0, 480.0, 415.0, 1.0),
    @(1.0, 0.0, 565.0, 575.0, 1.0),
    @(34.0, 0.0, 700.0, 580.0, 1.0),
    @(35.0, 0.0, 685.0, 595.0, 1.0),
    @(36.0, 0.0, 685.0, 610.0, 1.0),
    @(37.0, 0.0, 770.0, 610.0, 1.0),
    @(48.0, 0.0, 830.0, 610.0, 1.0),
    @(24.0, 0.0, 835.0, 625.0, 1.0),
    @(5.0, 0.0, 845.0, 655.0, 1.0),
    @(25.0, 0.0, 975.0, 580.0, 1.0),
    @(4.0, 0.0, 945.0, 685.0, 1.0),
    @(33.0, 0.0, 1150.0, 1160.0, 1.0),
    @(43.0, 0.0, 875.0, 920.0, 1.0),
    @(6.0, 0.0, 880.0, 660.0, 1.0),
    @(15.0, 0.0, 845.0, 680.0, 1.0),
    @(38.0, 0.0, 795.0, 645.0, 1.0),
    @(40.0, 0.0, 760.0, 650.0, 1.0),
    @(39.0, 0.0, 720.0, 635.0, 1.0),
    @(49.0, 0.0, 605.0, 625.0, 1.0),
    @(32.0, 0.0, 575.0, 665.0, 1.0),
    @(45.0, 0.0, 555.0, 815.0, 1.0),
    @(19.0, 0.0, 510.0, 875.0, 1.0),
    @(10.0, 0.0, 650.0, 1130.0, 1.0),
    @(9.0, 0.0, 580.0, 1175.0, 1.0),
    @(8.0, 0.0, 525.0, 1000.0, 1.0),
    @(41.0, 0.0, 475.0, 960.0, 1.0),
    @(22.0, 0.0, 520.0, 585.0, 1.0),
    @(31.0, 0.0, 420.0, 555.0, 1.0),
    @(18.0, 0.0, 415.0, 635.0, 1.0),
    @(3.0, 0.0, 345.0, 750.0, 1.0),
    @(21.0, 0.0, 300.0, 465.0, 1.0),
    @(17.0, 0.0, 145.0, 665.0, 1.0),
    @(42.0, 0.0, 95.0, 260.0, 1.0),
    @(7.0, 0.0, 25.0, 230.0, 1.0),
    @(2.0, 0.0, 25.0, 185.0, 1.0)
)

$r = 4
foreach ($row in $tour) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
